# Fruta / hortaliza, semanal
# Insert a new data row at row 338 (pushing the existing rows 338-381 down
# to 339-382) and populate it with the new observation.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 338..381 down to 339..382, leaving a blank row 338 to fill in.
$ws.Rows.Item(338).Insert()

$ws.Range("A338").Value = 3
$ws.Range("B338").Value = "Femacal de La Calera"
$ws.Range("C338").Value = "Coquimbo"
$ws.Range("D338").Value = 45212
$ws.Range("E338").Value = 5
$ws.Range("F338").Value = "Fruta"
$ws.Range("G338").Value = 100101
$ws.Range("H338").Value = "Berries"
$ws.Range("I338").Value = 100101001
$ws.Range("J338").Value = "Arándano (blue)"
$ws.Range("K338").Value = "Sin especificar"
$ws.Range("L338").Value = "Primera"
$ws.Range("M338").Value = 40
$ws.Range("N338").Value = 12000
$ws.Range("O338").Value = 12000
$ws.Range("P338").Value = 12000
$ws.Range("Q338").Value = "`$/bandeja 2 kilos"
$ws.Range("R338").Value = "Provincia de Quillota"
$ws.Range("S338").Value = 6000
$ws.Range("T338").Value = 2
